# Remove the "number_classes_2D" column from the Classes2D sheet and the
# "number_classes_3D" column from the Classes3D sheet. Both were the first
# (leftmost) column on their respective sheet, so deleting column A shifts
# the remaining columns left by one and shrinks the used range accordingly.

$wb = $excel.ActiveWorkbook

$classes2D = $wb.Worksheets.Item("Classes2D")
$classes2D.Columns.Item(1).Delete()

$classes3D = $wb.Worksheets.Item("Classes3D")
$classes3D.Columns.Item(1).Delete()
